# Weekly update for "Fruta, Vega Central Mapocho de Santiago - Granada"
# Insert two new rows of data (week of 2023-05-24, origin "Paine") right
# after the existing row 27, pushing the remaining historical rows down
# by two positions (old row 28 becomes 30, ..., old row 46 becomes 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 28 (existing row 27 stays put, the
# rest of the data shifts down by two rows).
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# Common values shared by every record in this sheet.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100104
$producto    = "Frutos de pepita"
$categoriaId = 100104001
$categoria   = "Granada"
$variedad    = "Wonderfull"
$fecha       = "2023-05-24"

# New row 28: Calidad "Especial"
$ws.Cells.Item(28, 1).Value  = $mercadoId
$ws.Cells.Item(28, 2).Value  = $mercado
$ws.Cells.Item(28, 3).Value  = $region
$ws.Cells.Item(28, 4).Value  = $fecha
$ws.Cells.Item(28, 5).Value  = $codreg
$ws.Cells.Item(28, 6).Value  = $tipo
$ws.Cells.Item(28, 7).Value  = $productoId
$ws.Cells.Item(28, 8).Value  = $producto
$ws.Cells.Item(28, 9).Value  = $categoriaId
$ws.Cells.Item(28, 10).Value = $categoria
$ws.Cells.Item(28, 11).Value = $variedad
$ws.Cells.Item(28, 12).Value = "Especial"
$ws.Cells.Item(28, 13).Value = 280
$ws.Cells.Item(28, 14).Value = 10500
$ws.Cells.Item(28, 15).Value = 10500
$ws.Cells.Item(28, 16).Value = 10500
$ws.Cells.Item(28, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(28, 18).Value = "Paine"
$ws.Cells.Item(28, 19).Value = 700
$ws.Cells.Item(28, 20).Value = 15

# New row 29: Calidad "Primera"
$ws.Cells.Item(29, 1).Value  = $mercadoId
$ws.Cells.Item(29, 2).Value  = $mercado
$ws.Cells.Item(29, 3).Value  = $region
$ws.Cells.Item(29, 4).Value  = $fecha
$ws.Cells.Item(29, 5).Value  = $codreg
$ws.Cells.Item(29, 6).Value  = $tipo
$ws.Cells.Item(29, 7).Value  = $productoId
$ws.Cells.Item(29, 8).Value  = $producto
$ws.Cells.Item(29, 9).Value  = $categoriaId
$ws.Cells.Item(29, 10).Value = $categoria
$ws.Cells.Item(29, 11).Value = $variedad
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 220
$ws.Cells.Item(29, 14).Value = 7500
$ws.Cells.Item(29, 15).Value = 7500
$ws.Cells.Item(29, 16).Value = 7500
$ws.Cells.Item(29, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(29, 18).Value = "Paine"
$ws.Cells.Item(29, 19).Value = 500
$ws.Cells.Item(29, 20).Value = 15
